$d = $word.ActiveDocument
$hf = $d.Sections.Item(1).Headers.Item(1)
$rng = $hf.Range
$f = $rng.Fields.Item(1)
$f.Delete()

$insPoint = $hf.Range.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertAfter("valueofx")

$textRange = $hf.Range.Duplicate
$textRange.Collapse(1)
$textRange.MoveEnd(1, 8)
Write-Output ("textRange text: [" + $textRange.Text + "]")
$textRange.Font.TextColor.ObjectThemeColor = 9
$textRange.Font.TextColor.RGB = 4626167
Write-Output "applied both via ColorFormat.RGB"
